$wb = $excel.ActiveWorkbook

# --- Worksheets involved ---
$summary  = $wb.Worksheets.Item("summary")
$baseline = $wb.Worksheets.Item("baseline")

# --- Fill in the newly-collected "6_wm" test row (row 9) on the summary sheet ---
# Apache Solr (D/E), Hammer (H/I), Google (L/M), Google CSE (P/Q) precision/recall,
# plus the execution-time columns (T/U). The dependent AVERAGE() formulas in row 13
# recalculate automatically.
$summary.Range("D9").Value = 1
$summary.Range("E9").Value = 1
$summary.Range("H9").Value = 1
$summary.Range("I9").Value = 0.037
$summary.Range("L9").Value = 1
$summary.Range("M9").Value = 1
$summary.Range("P9").Value = 1
$summary.Range("Q9").Value = 0.0515
$summary.Range("T9").Value = 335
$summary.Range("U9").Value = 1137

# --- View-state updates ---
# baseline sheet: move the selection, and drop its "selected tab" flag since the
# workbook will come back up on the summary sheet.
$baseline.Activate()
$baseline.Range("D10").Select()

# summary sheet: becomes the active sheet/tab, selection moves to U9, and the
# sheet is scrolled back to show column A (topLeftCell reset to default).
$summary.Activate()
$summary.Range("U9").Select()
